$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1752
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 501
$ws.Range("F6").Value = 1281
$ws.Range("F7").Value = 373
$ws.Range("F9").Value = 905
$ws.Range("F10").Value = 723
$ws.Range("F11").Value = 201
$ws.Range("F12").Value = 533
$ws.Range("F15").Value = 172
$ws.Range("F16").Value = 3012
$ws.Range("F17").Value = 2653
$ws.Range("F19").Value = 30
$ws.Range("F21").Value = 323
$ws.Range("F22").Value = 240
$ws.Range("F24").Value = 5390
$ws.Range("F28").Value = 61
$ws.Range("F29").Value = 354
$ws.Range("F30").Value = 1131
$ws.Range("F32").Value = 70
$ws.Range("F33").Value = 301

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1155
$ws.Range("F10").Value = 39
$ws.Range("F21").Value = 4
$ws.Range("F26").Value = 3990

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2512
$ws.Range("F9").Value = 1365

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 2512
$ws.Range("F6").Value = 1752
$ws.Range("F8").Value = 1365
$ws.Range("F11").Value = 171
$ws.Range("F12").Value = 501
$ws.Range("F13").Value = 1281
$ws.Range("F14").Value = 373
$ws.Range("F15").Value = 905
$ws.Range("F16").Value = 723
$ws.Range("F17").Value = 1155
$ws.Range("F18").Value = 201
$ws.Range("F19").Value = 533
$ws.Range("F22").Value = 3012
$ws.Range("F23").Value = 2653
$ws.Range("F24").Value = 30
$ws.Range("F25").Value = 323
$ws.Range("F26").Value = 39
$ws.Range("F27").Value = 240
$ws.Range("F29").Value = 5390
$ws.Range("F34").Value = 61
$ws.Range("F35").Value = 354
$ws.Range("F38").Value = 4
$ws.Range("F42").Value = 1131
$ws.Range("F47").Value = 70
$ws.Range("F48").Value = 301
